$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append new country rows (Belgium, Luxembourg, Canada) to the table
$ws.Range("A34").Value = "BE"
$ws.Range("B34").Value = "Belgium"

$ws.Range("A35").Value = "LU"
$ws.Range("B35").Value = "Luxembourg"

$ws.Range("C34").Value = "Belgium "
$ws.Range("C35").Value = "Luxembourg "

$ws.Range("A36").Value = "CA"
$ws.Range("B36").Value = "Canada"
$ws.Range("C36").Value = "Canada "

# Update the visible window / selection to match the saved view
$ws.Range("C36").Select()
$excel.ActiveWindow.ScrollRow = 11
